$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 190, shifting the
# existing rows 190-206 down to 192-208 (style of row above, including the
# date number format on column D, is carried along automatically).
$ws.Rows.Item(190).EntireRow.Insert()
$ws.Rows.Item(190).EntireRow.Insert()

# Populate new row 190.
$ws.Range("A190").Value = 11
$ws.Range("B190").Value = 'Vega Monumental Concepción'
$ws.Range("C190").Value = 'Bíobío'
$ws.Range("D190").Value = 45013
$ws.Range("E190").Value = 8
$ws.Range("F190").Value = 'Fruta'
$ws.Range("G190").Value = 100109
$ws.Range("H190").Value = 'Uva'
$ws.Range("I190").Value = 100109001
$ws.Range("J190").Value = 'Uva'
$ws.Range("K190").Value = 'Red Globe'
$ws.Range("L190").Value = 'Primera'
$ws.Range("M190").Value = 100
$ws.Range("N190").Value = 10000
$ws.Range("O190").Value = 11000
$ws.Range("P190").Value = 10500
$ws.Range("Q190").Value = '$/bandeja 18 kilos'
$ws.Range("R190").Value = 'Región Metropolitana'
$ws.Range("S190").Value = 583
$ws.Range("T190").Value = 18

# Populate new row 191.
$ws.Range("A191").Value = 11
$ws.Range("B191").Value = 'Vega Monumental Concepción'
$ws.Range("C191").Value = 'Bíobío'
$ws.Range("D191").Value = 45013
$ws.Range("E191").Value = 8
$ws.Range("F191").Value = 'Fruta'
$ws.Range("G191").Value = 100109
$ws.Range("H191").Value = 'Uva'
$ws.Range("I191").Value = 100109001
$ws.Range("J191").Value = 'Uva'
$ws.Range("K191").Value = 'Thompson seedless'
$ws.Range("L191").Value = 'Primera'
$ws.Range("M191").Value = 100
$ws.Range("N191").Value = 14000
$ws.Range("O191").Value = 15000
$ws.Range("P191").Value = 14500
$ws.Range("Q191").Value = '$/bandeja 18 kilos'
$ws.Range("R191").Value = 'Región Metropolitana'
$ws.Range("S191").Value = 806
$ws.Range("T191").Value = 18
